# Regenerate save_data to use K instead of Strike#: update column G ("K")
# values for the existing data rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "K" values (column G) keyed by row number, as recalculated.
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 2
    8  = 0
    9  = 1
    10 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
